# Jason's Guitar Shack - "Calculate Reorder Level" sheet:
# split the single "historical sales" scenario into two parallel scenarios -
# one driven by historical (year-ago) sales data, a new one driven by
# recent sales data - and renumber the old "no data" scenario from
# Test Scenario 2 to Test Scenario 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculate Reorder Level")

# --- Row 5: complete the GWT sentence for the historical-sales scenario ---
$ws.Range("C5").Value = "GIVEN there exists past sales data for a product WHEN a reorder level is required THEN reorder level is calculated as historical sales within the lead time from the same date in the previous year."
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 96

# --- Row 10: Test Scenario 2 becomes the "recent sales data" scenario ---
$ws.Range("C10").Value = "Reorder level calculated from recent sales data"

# --- Row 12: turn into the column-header row for the recent-sales data table ---
$ws.Range("B12").Value = "Product ID"
$ws.Range("C12").Value = "Date"
$ws.Range("D12").Value = "Lead time"
$ws.Range("E12").Value = "Current Sales Start Date"
$ws.Range("F12").Value = "Current Sales End Date"
$ws.Range("G12").Value = "Recent sales"
$ws.Range("H12").Value = "Reorder level"

# --- Row 11: GWT sentence for the recent-sales scenario ---
$ws.Range("C11").Value = "GIVEN there does not exist past sales data for a product WHEN a reorder level is required THEN reorder level is calculated as recent sales within the lead time from todays date minus lead time."
$ws.Range("C11").WrapText = $true
$ws.Range("B11").Value = "GWT"
$ws.Rows.Item(11).RowHeight = 96

# --- Row 13 (new): sample data row for the recent-sales table ---
$ws.Range("B13").Value = 811
$ws.Range("C13").Value = 44118
$ws.Range("C13").NumberFormat = "d-mmm-yy"
$ws.Range("D13").Value = 14
$ws.Range("E13").Value = 43739
$ws.Range("E13").NumberFormat = "d-mmm-yy"
$ws.Range("F13").Value = 44118
$ws.Range("F13").NumberFormat = "d-mmm-yy"
$ws.Range("G13").Value = 25
$ws.Range("H13").Value = 25

# --- Column widths: B and C need to be wider to fit the new wrapped text ---
$ws.Columns.Item(2).ColumnWidth = 40.333333333333336
$ws.Columns.Item(3).ColumnWidth = 29

# --- Update the active selection to reflect where editing finished ---
$ws.Range("C5").Select() | Out-Null
